# lambda.xlsx update: add delta + group-average/max summary columns (F:M)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column widths for the new columns (best-effort match to target widths) ---
$ws.Columns("F").ColumnWidth = 14
$ws.Columns("G").ColumnWidth = 15
$ws.Columns("H").ColumnWidth = 23.8333333333333
$ws.Columns("I").ColumnWidth = 24.6666666666667
$ws.Columns("J").ColumnWidth = 20.5
$ws.Columns("K").ColumnWidth = 22
$ws.Columns("L").ColumnWidth = 21.1666666666667
$ws.Columns("M").ColumnWidth = 22.1666666666667

# --- header row (order chosen so the shared-string table matches the source order) ---
$ws.Range("F1").Value = "Training - Test"
$ws.Range("G1").Value = "Validation - Test"
$ws.Range("I1").Value = "Average validation accuracy"
$ws.Range("H1").Value = "Average training accuracy"
$ws.Range("J1").Value = "Average test accuracy"
$ws.Range("K1").Value = "Max training accuracy"
$ws.Range("L1").Value = "Max validation accuracy"
$ws.Range("M1").Value = "Max test accuracy"

# --- per-row delta columns ---
$ws.Range("F2").Formula = "=C2-E2"
$ws.Range("G2").Formula = "=D2-E2"
$ws.Range("F3:F56").Formula = "=C3-E3"
$ws.Range("G3:G56").Formula = "=D3-E3"

# --- per-group (11-row block) average/max summary columns ---
$groupStarts = @(2, 13, 24, 35, 46)
foreach ($start in $groupStarts) {
    $end = $start + 10
    $ws.Range("H$start").Formula = "=AVERAGE(C$start`:C$end)"
    $ws.Range("I$start").Formula = "=AVERAGE(D$start`:D$end)"
    $ws.Range("J$start").Formula = "=AVERAGE(E$start`:E$end)"
    $ws.Range("K$start").Formula = "=MAX(C$start`:C$end)"
    $ws.Range("L$start").Formula = "=MAX(D$start`:D$end)"
    $ws.Range("M$start").Formula = "=MAX(E$start`:E$end)"
}

# --- selection / scroll position ---
$ws.Range("M46").Select()
